$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDLE")
$ws.Range("B2").Value = -2
